# refactor: Improve filter logic; remove V1 artifacts
#
# - Rename the "Customer" header (B1) to "Project" (the column's data -
#   Customer 1..4 - is kept as-is, only the label changes).
# - Remove the V1 "Non-Billable" flag column (old column F, which only ever
#   held a stray "X" marker or was blank) by selecting the whole column and
#   deleting it outright, shifting every later column (old G..K) one to the
#   left (new F..J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the Customer column header to Project.
$ws.Range("B1").Value = "Project"

# Select column F (the old "Non-Billable" marker column) and delete it
# entirely, shifting everything to its right one column to the left.
$ws.Range("F:F").Select()
$ws.Range("F:F").Delete()
